$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1/J1 - values first
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the existing header formatting (bold/border/center) by copying
# the format from the existing H1 header cell, same as the other headers.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# New data columns I and J, rows 2-8
$ws.Range("I2").Value = 6
$ws.Range("J2").Value = 9

$ws.Range("I3").Value = 6
$ws.Range("J3").Value = 9

$ws.Range("I4").Value = 2
$ws.Range("J4").Value = 5

$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 1

$ws.Range("I6").Value = 8
$ws.Range("J6").Value = 9

$ws.Range("I7").Value = 8
$ws.Range("J7").Value = 9

$ws.Range("I8").Value = 7
$ws.Range("J8").Value = 7
